# Weekly update: 3 new daily price records added for
# "Hortaliza, Vega Central Mapocho de Santiago - Puerro".
# Each new record is inserted as a new row (shifting the existing rows
# below it down by one), keeping the sheet's newest-first-ish ordering
# exactly as produced by the upstream daily/weekly consolidation job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-PuerroRow {
    param(
        [int]$RowIndex,
        [double]$Fecha,
        [double]$Volumen,
        [double]$PrecioMinimo,
        [double]$PrecioMaximo,
        [double]$PrecioPromedio,
        [double]$PrecioKg
    )

    # Push existing data (from $RowIndex downward) one row down, leaving
    # a blank row at $RowIndex for the new record.
    $ws.Rows.Item($RowIndex).Insert()

    $ws.Cells.Item($RowIndex, 1).Value = 9
    $ws.Cells.Item($RowIndex, 2).Value = 'Vega Central Mapocho de Santiago'
    $ws.Cells.Item($RowIndex, 3).Value = 'Metropolitana'
    $ws.Cells.Item($RowIndex, 4).Value = $Fecha
    $ws.Cells.Item($RowIndex, 5).Value = 13
    $ws.Cells.Item($RowIndex, 6).Value = 100112005
    $ws.Cells.Item($RowIndex, 7).Value = 'Puerro'
    $ws.Cells.Item($RowIndex, 8).Value = 'Sin especificar'
    $ws.Cells.Item($RowIndex, 9).Value = 'Primera'
    $ws.Cells.Item($RowIndex, 10).Value = $Volumen
    $ws.Cells.Item($RowIndex, 11).Value = $PrecioMinimo
    $ws.Cells.Item($RowIndex, 12).Value = $PrecioMaximo
    $ws.Cells.Item($RowIndex, 13).Value = $PrecioPromedio
    $ws.Cells.Item($RowIndex, 14).Value = '$/paquete 20 unidades'
    $ws.Cells.Item($RowIndex, 15).Value = 'Provincia de Chacabuco'
    $ws.Cells.Item($RowIndex, 16).Value = $PrecioKg
    $ws.Cells.Item($RowIndex, 17).Value = 20
    $ws.Cells.Item($RowIndex, 18).Value = 'Hortaliza'
}

# Insert from top to bottom so each target row index already accounts for
# the rows inserted earlier in this same pass.
Add-PuerroRow 12 44630 79  9000 10000 9494 475
Add-PuerroRow 50 44679 97  8000 9000  8505 425
Add-PuerroRow 84 44650 160 9000 10000 9500 475
